# chore(runtime): publish files + archive (2025-11-04 11:02:35)
#
# Adds three new KHL matches (2025-11-03) to Matches_SOG, and refreshes the
# derived Shots_HA / Shots_Summary aggregates + Meta_ext "as_of" snapshot to
# account for them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matches_SOG: append rows 427-429 for the 2025-11-03 games
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches_SOG")

# uid column (A) holds numbers-that-look-like-ids but is stored as text in
# this sheet, so enter them with a leading apostrophe to force text.
$matches.Range("A427").Value = "'897725"
$matches.Range("B427").Value = "2025-11-03T14:30:00"
$matches.Range("C427").Value = "Металлург Мг"
$matches.Range("D427").Value = "Драконы"
$matches.Range("E427").Value = 44
$matches.Range("F427").Value = 24
$matches.Range("G427").Value = "khl_text"

$matches.Range("A428").Value = "'897727"
$matches.Range("B428").Value = "2025-11-03T17:30:00"
$matches.Range("C428").Value = "Барыс"
$matches.Range("D428").Value = "Нефтехимик"
$matches.Range("E428").Value = 38
$matches.Range("F428").Value = 33
$matches.Range("G428").Value = "khl_text"

$matches.Range("A429").Value = "'897724"
$matches.Range("B429").Value = "2025-11-03T19:10:00"
$matches.Range("C429").Value = "Динамо Мн"
$matches.Range("D429").Value = "ЦСКА"
$matches.Range("E429").Value = 23
$matches.Range("F429").Value = 26
$matches.Range("G429").Value = "khl_text"

# ---------------------------------------------------------------------------
# 2) Shots_HA: bump as_of_utc on every team row, and refresh the home/away
#    on-goal totals for the six teams involved in the new games.
# ---------------------------------------------------------------------------
$shotsHA = $wb.Worksheets.Item("Shots_HA")
$newAsOf = "2025-11-03T19:10:00Z"

for ($r = 2; $r -le 23; $r++) {
    $shotsHA.Range("D" + $r).Value = $newAsOf
}

# Барыс (row 7)
$shotsHA.Range("E7").Value = 27
$shotsHA.Range("G7").Value = 850
$shotsHA.Range("H7").Value = 859
$shotsHA.Range("I7").Value = 31.5

# Динамо Мн (row 9)
$shotsHA.Range("E9").Value = 21
$shotsHA.Range("G9").Value = 753
$shotsHA.Range("H9").Value = 577
$shotsHA.Range("I9").Value = 35.9
$shotsHA.Range("J9").Value = 27.5

# Драконы (row 10)
$shotsHA.Range("F10").Value = 19
$shotsHA.Range("K10").Value = 527
$shotsHA.Range("L10").Value = 698
$shotsHA.Range("M10").Value = 27.7
$shotsHA.Range("N10").Value = 36.7

# Металлург Мг (row 13)
$shotsHA.Range("E13").Value = 25
$shotsHA.Range("G13").Value = 884
$shotsHA.Range("H13").Value = 646
$shotsHA.Range("I13").Value = 35.4
$shotsHA.Range("J13").Value = 25.8

# Нефтехимик (row 14)
$shotsHA.Range("F14").Value = 18
$shotsHA.Range("K14").Value = 498
$shotsHA.Range("L14").Value = 681
$shotsHA.Range("M14").Value = 27.7

# ЦСКА (row 23)
$shotsHA.Range("F23").Value = 21
$shotsHA.Range("K23").Value = 523
$shotsHA.Range("L23").Value = 597
$shotsHA.Range("N23").Value = 28.4

# ---------------------------------------------------------------------------
# 3) Shots_Summary: bump as_of_utc on every team row, and refresh the totals
#    for the same six teams.
# ---------------------------------------------------------------------------
$shotsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $shotsSummary.Range("D" + $r).Value = $newAsOf
}

# Барыс (row 7)
$shotsSummary.Range("E7").Value = 42
$shotsSummary.Range("F7").Value = 1265
$shotsSummary.Range("G7").Value = 1374
$shotsSummary.Range("H7").Value = 30.1

# Динамо Мн (row 9)
$shotsSummary.Range("E9").Value = 39
$shotsSummary.Range("F9").Value = 1421
$shotsSummary.Range("G9").Value = 1061
$shotsSummary.Range("H9").Value = 36.4

# Драконы (row 10)
$shotsSummary.Range("E10").Value = 38
$shotsSummary.Range("F10").Value = 1059
$shotsSummary.Range("G10").Value = 1360
$shotsSummary.Range("H10").Value = 27.9
$shotsSummary.Range("I10").Value = 35.8

# Металлург Мг (row 13)
$shotsSummary.Range("E13").Value = 39
$shotsSummary.Range("F13").Value = 1287
$shotsSummary.Range("G13").Value = 1009
$shotsSummary.Range("H13").Value = 33

# Нефтехимик (row 14)
$shotsSummary.Range("E14").Value = 41
$shotsSummary.Range("F14").Value = 1222
$shotsSummary.Range("G14").Value = 1451
$shotsSummary.Range("H14").Value = 29.8
$shotsSummary.Range("I14").Value = 35.4

# ЦСКА (row 23)
$shotsSummary.Range("E23").Value = 38
$shotsSummary.Range("F23").Value = 906
$shotsSummary.Range("G23").Value = 1096
$shotsSummary.Range("I23").Value = 28.8

# ---------------------------------------------------------------------------
# 4) Meta_ext: bump the snapshot timestamp and build_version counter.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta_ext")
$meta.Range("B2").Value = $newAsOf
$meta.Range("D2").Value = 39
